$wb = $excel.ActiveWorkbook

# CAP, CAP_NEW, and INVESTMENT sheets are rebuilt from scratch: the
# dimension shrinks from A1:P13 to A1:N13 (the 2015/2020 columns are
# dropped so the year header lines up with the REMOVAL sheet), and the
# DAC-scenario parameter values are reset to match (all zero).
$sheetNames = @("CAP", "CAP_NEW", "INVESTMENT")
$years = @(2025, 2030, 2035, 2040, 2045, 2050, 2055, 2060, 2070, 2080, 2090, 2100, 2110)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the old 2015/2020 columns (O/P after the shift) so the sheet
    # shrinks back down to A1:N13.
    $ws.Range("O1:P13").Delete()

    # Re-stamp the year header row (B1:N1) with the new year sequence.
    for ($i = 0; $i -lt $years.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item(1, $col).Value = $years[$i]
    }

    # Reset every data cell (rows 2-13) back to zero.
    $ws.Range("B2:N13").Value = 0
}

# REMOVAL keeps its existing A1:N13 layout/header, but its values are
# reset to zero to match the other recreated sheets.
$wsRemoval = $wb.Worksheets.Item("REMOVAL")
$wsRemoval.Range("B2:N13").Value = 0
